$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph and its index
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Docente\(s\) Responsável\(eis\)") {
        $targetIndex = $i
        break
    }
}

$target = $d.Paragraphs.Item($targetIndex)
$target.Range.InsertParagraphAfter()

# The newly inserted (empty) paragraph is now the next one
$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Range.Text = "5840535 - Messias Borges Silva"
$newPara.Style = "ListBullet"
